$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00721500721500722
$ws.Range("C2").Value = 0.543290043290043
$ws.Range("D2").Value = 0.0505050505050505
$ws.Range("E2").Value = 0.0555555555555556
$ws.Range("F2").Value = 0.00865800865800866
$ws.Range("G2").Value = 0.00793650793650794
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.0331890331890332
$ws.Range("J2").Value = 0.238095238095238
$ws.Range("K2").Value = 0.758297258297258
$ws.Range("L2").Value = 0.213564213564214
$ws.Range("M2").Value = 0.782828282828283
$ws.Range("N2").Value = 0.0108225108225108
$ws.Range("O2").Value = 0.0562770562770563
$ws.Range("P2").Value = 0.0158730158730159
$ws.Range("Q2").Value = 0.0281385281385281
$ws.Range("R2").Value = 0.000721500721500722
$ws.Range("S2").Value = 0.994227994227994
$ws.Range("T2").Value = 0.00216450216450216
$ws.Range("U2").Value = 0.209235209235209
$ws.Range("V2").Value = 0.00288600288600289
$ws.Range("W2").Value = 0.24963924963925
$ws.Range("X2").Value = 0.216450216450216
$ws.Range("B3").Value = 0.758297258297258
$ws.Range("C3").Value = 0.0202020202020202
$ws.Range("D3").Value = 0.018037518037518
$ws.Range("E3").Value = 0.218614718614719
$ws.Range("F3").Value = 0.210678210678211
$ws.Range("G3").Value = 0.992063492063492
$ws.Range("H3").Value = 0.784271284271284
$ws.Range("I3").Value = 0.963924963924964
$ws.Range("J3").Value = 0.00432900432900433
$ws.Range("K3").Value = 0.0230880230880231
$ws.Range("L3").Value = 0.00360750360750361
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0.000721500721500722
$ws.Range("O3").Value = 0.935786435786436
$ws.Range("P3").Value = 0.00144300144300144
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.968253968253968
$ws.Range("S3").Value = 0.000721500721500722
$ws.Range("T3").Value = 0.0238095238095238
$ws.Range("U3").Value = 0.788600288600289
$ws.Range("V3").Value = 0.784992784992785
$ws.Range("W3").Value = 0.00865800865800866
$ws.Range("X3").Value = 0.0108225108225108
$ws.Range("B4").Value = 0.207070707070707
$ws.Range("C4").Value = 0.427128427128427
$ws.Range("D4").Value = 0.20995670995671
$ws.Range("E4").Value = 0.713564213564214
$ws.Range("F4").Value = 0.776334776334776
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.00360750360750361
$ws.Range("I4").Value = 0.000721500721500722
$ws.Range("J4").Value = 0.756854256854257
$ws.Range("K4").Value = 0.215007215007215
$ws.Range("L4").Value = 0.779220779220779
$ws.Range("M4").Value = 0.00216450216450216
$ws.Range("N4").Value = 0.963924963924964
$ws.Range("O4").Value = 0.000721500721500722
$ws.Range("P4").Value = 0.000721500721500722
$ws.Range("Q4").Value = 0.971861471861472
$ws.Range("R4").Value = 0.00360750360750361
$ws.Range("S4").Value = 0.00432900432900433
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0.00216450216450216
$ws.Range("V4").Value = 0.000721500721500722
$ws.Range("W4").Value = 0.73015873015873
$ws.Range("X4").Value = 0.769119769119769
$ws.Range("B5").Value = 0.0274170274170274
$ws.Range("C5").Value = 0.00865800865800866
$ws.Range("D5").Value = 0.721500721500722
$ws.Range("E5").Value = 0.0122655122655123
$ws.Range("F5").Value = 0.00432900432900433
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.212121212121212
$ws.Range("I5").Value = 0.00216450216450216
$ws.Range("J5").Value = 0.000721500721500722
$ws.Range("K5").Value = 0.00360750360750361
$ws.Range("L5").Value = 0.00360750360750361
$ws.Range("M5").Value = 0.215007215007215
$ws.Range("N5").Value = 0.0245310245310245
$ws.Range("O5").Value = 0.00649350649350649
$ws.Range("P5").Value = 0.981962481962482
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0.0274170274170274
$ws.Range("S5").Value = 0.000721500721500722
$ws.Range("T5").Value = 0.974025974025974
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0.211399711399711
$ws.Range("W5").Value = 0.0115440115440115
$ws.Range("X5").Value = 0.00360750360750361
